# feat: add 2022-Q1 data
#
# - Insert a new sheet "2022-Q1" (fund-holding detail) right before the
#   "总计" (totals) summary sheet.
# - Prepend a "2022-Q1" row to the "总计" summary table.

$wb = $excel.ActiveWorkbook

$templateSheet = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# helper: write a value into a cell as literal TEXT (no auto number
# coercion, no extra style) by routing it through a scratch formula
# cell + PasteSpecial(xlPasteValues). Mirrors how the source data (an
# inline/shared string even though it "looks like" a number) is stored.
# ---------------------------------------------------------------------
function Set-TextValue($sheet, $cellRef, $text) {
    $scratch = $sheet.Range("ZZ1")
    $escaped = $text -replace '"', '""'
    $scratch.Formula = "=""$escaped"""
    $scratch.Copy()
    $sheet.Range($cellRef).PasteSpecial(-4163)
    $scratch.Clear()
}

# =======================================================================
# 1) New sheet "2022-Q1", inserted right before "总计"
# =======================================================================
# NOTE: inserting a sheet "Before:" an existing sheet reference re-points
# that *positional* handle at the freshly-inserted sheet, so the old
# "总计" worksheet object must be re-fetched by name afterwards (see
# part 2 below) rather than reused.
$totalSheetBeforeInsert = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheetBeforeInsert)
$newSheet.Name = "2022-Q1"

# Pull header (row1, B:H) and index-column (A2:A7) formatting from an
# existing fund-detail sheet so the new sheet matches the house style.
$templateSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$templateSheet.Range("A2:A7").Copy()
$newSheet.Range("A2:A7").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$fundRows = @(
    @{ Row=2; Code="519170"; Name="浦银安盛增长动力灵活配置混合"; Size="8.12"; Position="85.61"; Ratio="3.14"; Value="0.2550"; Rank=3  },
    @{ Row=3; Code="009169"; Name="湘财长兴灵活配置混合A";       Size="1.16"; Position="85.40"; Ratio="6.01"; Value="0.0697"; Rank=1  },
    @{ Row=4; Code="009170"; Name="湘财长兴灵活配置混合C";       Size="0.46"; Position="85.40"; Ratio="6.01"; Value="0.0276"; Rank=1  },
    @{ Row=5; Code="001648"; Name="工银瑞信新价值灵活配置混合"; Size="1.96"; Position="80.70"; Ratio="1.27"; Value="0.0249"; Rank=10 },
    @{ Row=6; Code="005537"; Name="中航新起航灵活配置混合A";     Size="0.03"; Position="87.09"; Ratio="8.44"; Value="0.0025"; Rank=2  },
    @{ Row=7; Code="005538"; Name="中航新起航灵活配置混合C";     Size="0.01"; Position="87.09"; Ratio="8.44"; Value="0.0008"; Rank=2  }
)

foreach ($r in $fundRows) {
    $row = $r.Row
    $newSheet.Range("A$row").Value = ($row - 2)
    Set-TextValue $newSheet "B$row" $r.Code
    $newSheet.Range("C$row").Value = $r.Name
    Set-TextValue $newSheet "D$row" $r.Size
    Set-TextValue $newSheet "E$row" $r.Position
    Set-TextValue $newSheet "F$row" $r.Ratio
    Set-TextValue $newSheet "G$row" $r.Value
    $newSheet.Range("H$row").Value = $r.Rank
}

# =======================================================================
# 2) Prepend a "2022-Q1" row to the "总计" summary sheet, shifting the
#    existing quarters down by one row.
# =======================================================================
$b4 = $totalSheet.Range("B4").Value()
$c4 = $totalSheet.Range("C4").Value()
$d4 = $totalSheet.Range("D4").Value()

$b3 = $totalSheet.Range("B3").Value()
$c3 = $totalSheet.Range("C3").Value()
$d3 = $totalSheet.Range("D3").Value()

$b2 = $totalSheet.Range("B2").Value()
$c2 = $totalSheet.Range("C2").Value()
$d2 = $totalSheet.Range("D2").Value()

# Row 5 (new) <- old row 4 data, copy formatting of A5 from A4 first.
$totalSheet.Range("A4").Copy()
$totalSheet.Range("A5").PasteSpecial(-4122)
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = $b4
$totalSheet.Range("C5").Value = $c4
$totalSheet.Range("D5").Value = $d4

# Row 4 <- old row 3 data.
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = $b3
$totalSheet.Range("C4").Value = $c3
$totalSheet.Range("D4").Value = $d3

# Row 3 <- old row 2 data.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = $b2
$totalSheet.Range("C3").Value = $c2
$totalSheet.Range("D3").Value = $d2

# Row 2 <- brand-new 2022-Q1 entry.
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 6
$totalSheet.Range("D2").Value = 0.38
